# --- "General" sheet: restructure from 3-col (idx,label,value) to 2-col (label,value) ---
$wb = $excel.ActiveWorkbook
$general = $wb.Worksheets.Item("General")

$general.Cells.Clear()

$general.Range("A1").Value = "Excel file version"
$general.Range("B1").Value = "v2.0"

$general.Range("A2").Value = "Name"
$general.Range("B2").Value = "ieee13pv"

$general.Range("A3").Value = "Frequency (Hz)"
$general.Range("B3").Value = 60

$general.Range("A4").Value = "Power Base (MVA)"
$general.Range("B4").Value = 100

# --- Add new "Switch" sheet right after "General" ---
$switch = $wb.Worksheets.Add($null, $general, 1, $null)
$switch.Name = "Switch"

# Header row (bold, thin border all sides, centered horizontally, top-aligned vertically)
$switch.Range("A1").Value = "From Bus"
$switch.Range("B1").Value = "To Bus"
$switch.Range("C1").Value = "ID"
$switch.Range("D1").Value = "Status"

$header = $switch.Range("A1:D1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Data rows — bus/device names that look like plain numbers ("650", "671692", ...) must
# stay text, so use a leading apostrophe (Excel's own "force text" convention) for those;
# alphanumeric identifiers are left as plain text assignments. ClearFormats() afterwards
# drops the quote-prefix formatting flag left behind, so no stray style sticks to the cell.
$switch.Range("A2").Value = "'650"
$switch.Range("B2").Value = "brkr"
$switch.Range("C2").Value = "brkr1"
$switch.Range("D2").Value = 1

$switch.Range("A3").Value = "'633"
$switch.Range("B3").Value = "xf1"
$switch.Range("C3").Value = "fuse1"
$switch.Range("D3").Value = 1

$switch.Range("A4").Value = "'671"
$switch.Range("B4").Value = "'692"
$switch.Range("C4").Value = "'671692"
$switch.Range("D4").Value = 1

$switch.Range("A5").Value = "'684"
$switch.Range("B5").Value = "tap"
$switch.Range("C5").Value = "sect1"
$switch.Range("D5").Value = 1

$switch.Range("A6").Value = "'632"
$switch.Range("B6").Value = "mid"
$switch.Range("C6").Value = "rec1"
$switch.Range("D6").Value = 1

$switch.Range("A2:D6").ClearFormats()

$general.Activate()
